$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10: zeroshot / full random / sentences / 3000 / 200 / 5 / (no G) / 4700 / null / yes / 110
$ws.Range("A10").Value = "zeroshot"
$ws.Range("B10").Value = "full random"
$ws.Range("C10").Value = "sentences"
$ws.Range("D10").Value = 3000
$ws.Range("E10").Value = 200
$ws.Range("F10").Value = 5
$ws.Range("H10").Value = 4700
$ws.Range("I10").Value = "null"
$ws.Range("J10").Value = "yes"
$ws.Range("K10").Value = 110

# Row 11: zeroshot huang / full random / sentences / 3000 / 200 / 5 / (no H) / null / yes / (no K)
$ws.Range("A11").Value = "zeroshot huang"
$ws.Range("B11").Value = "full random"
$ws.Range("C11").Value = "sentences"
$ws.Range("D11").Value = 3000
$ws.Range("E11").Value = 200
$ws.Range("F11").Value = 5
$ws.Range("I11").Value = "null"
$ws.Range("J11").Value = "yes"

# Update selection to match target (J13)
$ws.Range("J13").Select()
